$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1651255589585369"
$ws1.Range("B2").Value = "go_stims-1651255589554117.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555895697489.csv"
$ws1.Range("B4").Value = "go_stims-16512555895697489.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255589585369.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16512555928877506"
$ws2.Range("B2").Value = "TB-16512555920235386.csv"
$ws2.Range("B3").Value = "TB-16512555913047235.csv"
$ws2.Range("B4").Value = "ZB-match_6-16512555900254197.csv"
$ws2.Range("B5").Value = "OB-16512555912874486.csv"
$ws2.Range("B6").Value = "ZB-match_1-1651255589663489.csv"
$ws2.Range("B7").Value = "ZB-match_3-1651255589854782.csv"
$ws2.Range("B8").Value = "OB-16512555902786858.csv"
$ws2.Range("B9").Value = "OB-16512555906413271.csv"
$ws2.Range("B10").Value = "TB-1651255592872124.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16512555928877506"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512555929346275"
$ws4.Range("B2").Value = "MM_stims-16512555929033782.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555928877506.csv"
$ws4.Range("B4").Value = "MM_stims-16512555929190018.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555929033782.csv"
$ws4.Range("B6").Value = "MM_stims-16512555929346275.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555929190018.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512555929971282"
$ws5.Range("B2").Value = "SAT_stims-16512555929346275.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555929815106.csv"
$ws5.Range("B4").Value = "SAT_stims-1651255592950252.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555929658785.csv"
